$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" text on sheet "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $wsHoja1.Range("A1")
$text = $cell.Value2
$text = $text -replace [regex]::Escape("1000 Bs = 7.14 = 28808.17 pesos"), "1000 Bs = 6.98 = 28051.25 pesos"
$text = $text -replace [regex]::Escape("28808.17 pesos = 7.12 = 965.73 Bs"), "28051.25 pesos = 6.93 = 963.23 Bs"
$cell.Value2 = $text

# --- Update the rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 143.22
$wsTasas.Range("O10").Value = 4017.5
$wsTasas.Range("N12").Value = 4047.97
$wsTasas.Range("O12").Value = 139
